$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 0

$ws.Range("B3").Value = 0.770173668851462
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = -0.6961419073433766

$ws.Range("B4").Value = 0
$ws.Range("C4").Value = 0
$ws.Range("D4").Value = -0.7579308124743087

$ws.Range("B5").Value = 0
$ws.Range("C5").Value = -0.6950419217015434
$ws.Range("D5").Value = 0

$ws.Range("B6").Value = 0
$ws.Range("C6").Value = 0
$ws.Range("D6").Value = 0

$ws.Range("B7").Value = -0.8830107338985754
$ws.Range("C7").Value = 0
$ws.Range("D7").Value = 0.77232977548935

$ws.Range("B8").Value = 0
$ws.Range("C8").Value = 0
$ws.Range("D8").Value = -0.6627523733933205

$ws.Range("B9").Value = 0
$ws.Range("C9").Value = -0.6647337829335498
$ws.Range("D9").Value = 0
